$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.809.41"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "3.494.86"
$ws.Range("E3").Value = "  -2.17%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'607.62"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").Value = "'195.13"
$ws.Range("E6").Value = "  +3.49%  "
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.211"
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("E10").Value = "  +1.21%  "
$ws.Range("D11").Value = "'53.60"
$ws.Range("E11").Value = "  -1.49%  "
$ws.Range("D12").Value = "'0.0000307"
$ws.Range("E12").Value = "  -1.73%  "
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").Value = "4.058.92"
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("D15").Value = "'601.30"
$ws.Range("E15").Value = "  +4.44%  "
$ws.Range("D16").Value = "69.893.61"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("B17").Value = "Chainlink"
$ws.Range("C17").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D17").Value = "'18.93"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").Value = "'12.65"
$ws.Range("E18").Value = "  -1.44%  "
$ws.Range("D19").Value = "3.494.25"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("E21").Value = "  -1.12%  "
$ws.Range("D22").Value = "'17.97"
$ws.Range("E22").Value = "  +2.19%  "
$ws.Range("E23").Value = "  +11.00%  "
$ws.Range("D24").Value = "'4.65"
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("D25").Value = "'5.07"
$ws.Range("E25").Value = "  +3.68%  "
$ws.Range("E26").Value = "  +4.24%  "
$ws.Range("D27").Value = "'10.94"
$ws.Range("E27").Value = "  -0.44%  "
$ws.Range("D28").Value = "'9.80"
$ws.Range("E28").Value = "  +3.85%  "
$ws.Range("D29").Value = "'33.96"
$ws.Range("E29").Value = "  +4.61%  "
$ws.Range("D30").Value = "'4.44"
$ws.Range("E30").Value = "  +18.63%  "
$ws.Range("D31").Value = "'7.19"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").Value = "'12.67"
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("E33").Value = "  +0.27%  "
$ws.Range("D34").Value = "'64.08"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").Value = "3.685.68"
$ws.Range("E35").Value = "  -4.31%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  -0.18%  "
$ws.Range("D37").Value = "'519.24"
$ws.Range("E37").Value = "  -0.74%  "
$ws.Range("E38").Value = "  -5.60%  "
$ws.Range("D39").Value = "0.0₃0789"
$ws.Range("E39").Value = "  -0.16%  "
$ws.Range("D40").Value = "'0.391"
$ws.Range("E40").Value = "  -4.30%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("D43").Value = "'0.137"
$ws.Range("E43").Value = "  -1.64%  "
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("E45").Value = "  -4.34%  "
$ws.Range("E47").Value = "  -3.95%  "
$ws.Range("D48").Value = "'8.76"
$ws.Range("E48").Value = "  -5.13%  "
$ws.Range("E49").Value = "  +0.26%  "
$ws.Range("D50").Value = "'132.70"
$ws.Range("E50").Value = "  -2.88%  "
$ws.Range("D51").Value = "'1.29"
$ws.Range("E51").Value = "  +10.03%  "
